$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 'wrestling pants men'
$ws.Range("A2").Value = 'mens sliding pants'
$ws.Range("A3").Value = 'knee pads for baseball'
$ws.Range("A4").Value = 'sliding pad for baseball'
$ws.Range("A5").Value = 'compression pants padded knees'
$ws.Range("A6").Value = 'goalkeeper knee pads'
$ws.Range("A7").Value = 'tights for basketball youth'
$ws.Range("A8").Value = 'honeycomb basketball knee pads'
$ws.Range("A9").Value = 'knee pad shorts'
$ws.Range("A10").Value = 'mens compression capri leggings'
$ws.Range("A11").Value = 'boys compression pants with knee pads'
$ws.Range("A12").Value = 'rodilleras basketball'
$ws.Range("A13").Value = 'honeycomb knee pads basketball'
$ws.Range("A14").Value = 'pad pants'
$ws.Range("A15").Value = 'long basketball knee pads'
$ws.Range("A16").Value = 'men basketball knee'
$ws.Range("A17").Value = 'mens baseball sliding short'
$ws.Range("A18").Value = 'knee pad for basketball'
$ws.Range("A19").Value = 'hex knee pads'
$ws.Range("A20").Value = 'construction pants with knee pads'
$ws.Range("A21").Value = 'men volleyball knee pads'
$ws.Range("A22").Value = 'tight for boys basketball'
$ws.Range("A23").Value = 'bjj kneepads'
$ws.Range("A24").Value = 'boys padded compression pants'
$ws.Range("A25").Value = 'basketball knee tights'
$ws.Range("A26").Value = 'youth basketball compression leggings with knee pads'
$ws.Range("A27").Value = 'basketball knee compression'
$ws.Range("A28").Value = 'youth basketball compression tights'
$ws.Range("A29").Value = 'basketball hex knee pads'
$ws.Range("A30").Value = 'compression basketball leggings'
$ws.Range("A31").Value = 'knee pads for softball'
$ws.Range("A32").Value = 'the best knee pads'
$ws.Range("A33").Value = 'youth basketball knee pad'
$ws.Range("A34").Value = 'mens yoga pants compression'
$ws.Range("A35").Value = 'black knee pads basketball'
$ws.Range("A36").Value = 'boys basketball knee pads youth'
$ws.Range("A37").Value = 'tights basketball'
$ws.Range("A38").Value = 'men capri tights'
$ws.Range("A39").Value = 'nike basketball knee pads'
$ws.Range("A40").Value = 'padded yoga pants'
$ws.Range("A41").Value = 'padded knee pants'
$ws.Range("A42").Value = 'extra padded knee pads'
$ws.Range("A43").Value = 'basketball leg tights'
$ws.Range("A44").Value = 'catchers knee pads'
$ws.Range("A45").Value = 'basketball tights for youth boys'
$ws.Range("A46").Value = 'compression shorts with padding basketball'
$ws.Range("A47").Value = 'wrestling pants youth'
$ws.Range("A48").Value = 'honeycomb knee pads'
$ws.Range("A49").Value = 'baseball leggings for men'
$ws.Range("A50").Value = 'pants men basketball'
$ws.Range("A51").Value = 'long basketball knee pads adult'
$ws.Range("A52").Value = 'yellow leggings for men'
$ws.Range("A53").Value = 'mens compression capri'
$ws.Range("A54").Value = 'wrestling tights youth boy'
$ws.Range("A55").Value = 'basketball youth tights'
$ws.Range("A56").Value = 'need pads for basketball'
$ws.Range("A57").Value = 'compression capris for men'
$ws.Range("A58").Value = 'knee pads for youth basketball'
$ws.Range("A59").Value = 'knee pads under pants'
$ws.Range("A60").Value = 'goalkeeper padded pants'
$ws.Range("A61").Value = 'baseball knee pants'
$ws.Range("A62").Value = 'girls softball leggings'
$ws.Range("A63").Value = 'youth basketball tights'
$ws.Range("A64").Value = 'sliding pad'
$ws.Range("A65").Value = 'work pants with knee pads for men'
$ws.Range("A66").Value = 'multicam combat pants with knee pads'
$ws.Range("A67").Value = 'knee pad honeycomb'
$ws.Range("A68").Value = 'cycling leggings men'
$ws.Range("A69").Value = 'boys basketball compression tights'
$ws.Range("A70").Value = 'compression knee pads work'
$ws.Range("A71").Value = 'softball leggings'
$ws.Range("A72").Value = 'basketball compression tights'
$ws.Range("A73").Value = 'boys basketball pads'
$ws.Range("A74").Value = 'basketball knee pads youth boys black'
$ws.Range("A75").Value = 'lacrosse knee pads'
$ws.Range("A76").Value = 'knee pads youth basketball'
$ws.Range("A77").Value = 'hex kneepads'
$ws.Range("A78").Value = 'knee basketball'
$ws.Range("A79").Value = 'basketball compression pants'
$ws.Range("A80").Value = 'knee pad basketball youth'
$ws.Range("A81").Value = 'weightlifting pants men'
$ws.Range("A82").Value = 'mens knee baseball pants'
$ws.Range("A83").Value = 'thick yoga pad'
$ws.Range("A84").Value = 'padded basketball compression shorts'
$ws.Range("A85").Value = 'under pant knee pads'
$ws.Range("A86").Value = 'youth boys basketball knee pads'
$ws.Range("A87").Value = 'basketball pad'
$ws.Range("A88").Value = 'basketball tights for boys youth'
$ws.Range("A89").Value = 'ready man'
$ws.Range("A90").Value = 'youth basketball knee pads small'
$ws.Range("A91").Value = 'padded compression pants men'
$ws.Range("A92").Value = 'long sliding shorts baseball'
$ws.Range("A93").Value = 'basketball padded compression shorts men'
$ws.Range("A94").Value = 'compression capris men'
$ws.Range("A95").Value = 'gym men leggings'
$ws.Range("A96").Value = 'sliding shorts youth softball'
$ws.Range("A97").Value = 'mens compression pants basketball'
$ws.Range("A98").Value = 'softball tights'
$ws.Range("A99").Value = 'padded basketball shorts'
$ws.Range("A100").Value = 'knee pads basketball men'
